# Apply the edits described in the commit:
# "Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab from SCD0275 to SCD0017
$ws.Name = "SCD0017"

# Update the TC_ID values in column B (rows 2 and 3) from DGS-290 to SCD0017-005
$ws.Range("B2").Value = "SCD0017-005"
$ws.Range("B3").Value = "SCD0017-005"

# Widen column B to fit the new, longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.6667

# Update the view: zoom level and active selection cell
$ws.Range("B4").Select()
$excel.ActiveWindow.Zoom = 84
